$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three-language title row (A1:C1): drop the extra ".1" level
# from the indicator code "4.c.1.1" -> "4.c.1"
$ws.Range("A1").Value = "4.с.1 Билим берүү мекемелерде диплом берилгем мугалимдердин үлүшү"
$ws.Range("B1").Value = "4.c.1 Доля дипломированных учителей в образовательных учереждениях"
$ws.Range("C1").Value = "4.c.1 Proportion of certified teachers in educational institutions"

# Move the active cell selection from N11 to C11
$ws.Range("C11").Select()
